$d = $word.ActiveDocument

# The document currently ends with an empty paragraph (w:p/). We replace
# that paragraph's content with a new bold paragraph describing version
# management, matching the target OOXML exactly (including the
# firstLineChars/firstLine indent and the bCs run property, plus a
# trailing "." run carrying the eastAsia font hint / lang).

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$target = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes" ?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:ind w:firstLineChars="200" w:firstLine="442"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t xml:space="preserve">Version management is to manage the collection of specific functions or the construction results of specific codes in the process of software development, which mainly includes the management of version number, the preliminary planning of version, the response to the change of requirements during version development, and the summary and review after the release of version. Before version development: by establishing version number identification, clarifying version target, formulating version online requirements, and designing release strategy, product functions and quality can conform to user expectations as much as possible</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:b/>
                <w:bCs/>
                <w:lang w:eastAsia="zh-CN"/>
              </w:rPr>
              <w:t>.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$target.InsertXML($xml)

# InsertXML inserts the new paragraph before the paragraph-mark range that
# was addressed, leaving the original (now-empty) trailing paragraph
# behind it. Remove that leftover empty paragraph mark so the new
# paragraph becomes the final paragraph of the body, exactly as in the
# target document.
$trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
$cleanup = $d.Range($trailing.Range.Start - 1, $trailing.Range.End)
$cleanup.Delete()
